# Applies the diff: adds 3 new list items after "Priority Number", keeps the
# existing blank paragraph, then adds a block of narrative paragraphs, and
# finally appends one more blank paragraph after the bookmark paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the three new ListParagraph (ilvl 0, numId 1) bullet items
#    right after "Priority Number " (paragraph 13) and before the blank
#    paragraph that follows it.
# ---------------------------------------------------------------------
$priorityPara = $d.Paragraphs.Item(13)
$anchor = $priorityPara.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item(14)
$p1.Range.ListFormat.ListLevelNumber = 1
$p1.Range.Text = "Python for process data"

$anchor2 = $p1.Range
$anchor2.Collapse(0)
$anchor2.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(15)
$p2.Range.ListFormat.ListLevelNumber = 1
$p2.Range.Text = "ISS/express as webserver "

$anchor3 = $p2.Range
$anchor3.Collapse(0)
$anchor3.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(16)
$p3.Range.ListFormat.ListLevelNumber = 1
$p3.Range.Text = "Node.js"

# ---------------------------------------------------------------------
# 2) The next paragraph (17) is the pre-existing blank <w:p/>. Leave it
#    untouched, then insert the narrative block after it.
# ---------------------------------------------------------------------
$blankPara = $d.Paragraphs.Item(17)

$a = $blankPara.Range
$a.Collapse(0)
$a.InsertParagraphAfter()
$para1 = $d.Paragraphs.Item(18)
$para1.Range.ParagraphFormat.SpaceAfter = 0
$para1.Range.Text = "Data processing from excel to json injectable data for Node.js is Python "

$a = $para1.Range
$a.Collapse(0)
$a.InsertParagraphAfter()
$para2 = $d.Paragraphs.Item(19)
$para2.Range.ParagraphFormat.SpaceAfter = 0

$a = $para2.Range
$a.Collapse(0)
$a.InsertParagraphAfter()
$para3 = $d.Paragraphs.Item(20)
$para3.Range.ParagraphFormat.SpaceAfter = 0
$para3.Range.Text = "The json data is used by Node.js and create the page for express to display the page on the browser."

$a = $para3.Range
$a.Collapse(0)
$a.InsertParagraphAfter()
$para4 = $d.Paragraphs.Item(21)
$para4.Range.ParagraphFormat.SpaceAfter = 0

$a = $para4.Range
$a.Collapse(0)
$a.InsertParagraphAfter()
$para5 = $d.Paragraphs.Item(22)
$para5.Range.Text = "Excel -> python -> json -> Node.js -> HTML/CSS -> browser"

# ---------------------------------------------------------------------
# 3) Append a trailing blank paragraph after the bookmark paragraph
#    (which is now paragraph 23).
# ---------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item(23)
$a = $bookmarkPara.Range
$a.Collapse(0)
$a.InsertParagraphAfter()

Write-Output "done"
